# Update the "Date Placeholder" field shown on the slide masters / layouts /
# handout master / notes master from 1/21/2024 to 4/28/2025.
#
# ppPlaceholderDate = 16
$ppPlaceholderDate = 16
$NewDate = "4/28/2025"

function Update-DateShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.HasTextFrame) {
            if ($shp.PlaceholderFormat.Type -eq $ppPlaceholderDate) {
                $shp.TextFrame.TextRange.Text = $NewDate
            }
        }
    }
}

$p = $ppt.ActivePresentation

# 1. The slide master itself.
Update-DateShapes $p.SlideMaster.Shapes

# 2. Every slide layout ("custom layout") that hangs off the slide master.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShapes $layouts.Item($j).Shapes
}

# 3. The handout master.
Update-DateShapes $p.HandoutMaster.Shapes

# 4. The notes master.
Update-DateShapes $p.NotesMaster.Shapes
